# Applies the NEVADA_2024 cleaning edits:
#  1. Rename header columns to snake_case machine-readable names.
#  2. Title-case the Spanish connector words (de, del, la, las, el, los, y)
#     inside state/municipality names, and normalize "TOTAL" -> "Total".
#  3. Remove the trailing metadata/footer rows (1435-1439), shrinking the
#     used range down to A1:D1433.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-Text($s) {
    if ($s -eq "TOTAL") {
        return "Total"
    }
    $connectors = @("de", "del", "la", "las", "el", "los", "y")
    $words = $s -split " "
    $newWords = @()
    foreach ($w in $words) {
        if ($connectors -contains $w) {
            $newWords += ($w.Substring(0, 1).ToUpper() + $w.Substring(1))
        } else {
            $newWords += $w
        }
    }
    return ($newWords -join " ")
}

# --- 1. Header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- 2. Title-case state (col A) / municipality (col B) names ----------
for ($r = 2; $r -le 1433; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value = Transform-Text($a)
    }

    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null) {
        $ws.Cells.Item($r, 2).Value = Transform-Text($b)
    }
}

# --- 3. Drop the trailing metadata/footer rows --------------------------
$ws.Range("A1435:D1439").ClearContents()
